$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 09:52 AM"

# --- "1 Month Performance" sheet: refresh stock/% change rankings ---
$perf = $wb.Worksheets.Item("1 Month Performance")

$perf.Range("C4").Value = 77.1001
$perf.Range("C6").Value = 66.4884
$perf.Range("C7").Value = 65.36409999999999
$perf.Range("C8").Value = 64.92489999999999
$perf.Range("C10").Value = 51.0008
$perf.Range("C12").Value = 45.5488
$perf.Range("B13").Value = "MTARTECH"
$perf.Range("C13").Value = 40.8349
$perf.Range("B14").Value = "TVSSRICHAK"
$perf.Range("C14").Value = 40.5681
$perf.Range("C15").Value = 38.433
$perf.Range("C16").Value = 38.1142
$perf.Range("C17").Value = 36.8875
$perf.Range("C18").Value = 36.8344
$perf.Range("C19").Value = 36.47
$perf.Range("C20").Value = 36.3974
$perf.Range("B21").Value = "SOUTHBANK"
$perf.Range("C21").Value = 35.7662
$perf.Range("B22").Value = "ONMOBILE"
$perf.Range("C22").Value = 35.6918
$perf.Range("C23").Value = 35.6012
$perf.Range("B24").Value = "TVSELECT"
$perf.Range("C24").Value = 34.9546
$perf.Range("B25").Value = "MAANALU"
$perf.Range("C25").Value = 34.8477
$perf.Range("B26").Value = "RAMCOSYS"
$perf.Range("C26").Value = 34.7982
$perf.Range("B27").Value = "SHAREINDIA"
$perf.Range("C27").Value = 34.6917
$perf.Range("C30").Value = 31.4751
$perf.Range("C31").Value = 29.0155
$perf.Range("C32").Value = 28.9429
$perf.Range("B34").Value = "MINDTECK"
$perf.Range("C34").Value = 27.8676
$perf.Range("B35").Value = "CARTRADE"
$perf.Range("C35").Value = 27.8648
$perf.Range("C36").Value = 27.4033
$perf.Range("B37").Value = "HATSUN"
$perf.Range("C37").Value = 26.503
$perf.Range("B38").Value = "INDORAMA"
$perf.Range("C38").Value = 26.47
$perf.Range("C39").Value = 26.4498
$perf.Range("B40").Value = "IFBIND"
$perf.Range("C40").Value = 26.0965
$perf.Range("B41").Value = "MRPL"
$perf.Range("C41").Value = 26.0441
$perf.Range("C42").Value = 25.8791
$perf.Range("C43").Value = 25.7143
$perf.Range("C44").Value = 25.459
$perf.Range("C45").Value = 24.3704
$perf.Range("C46").Value = 24.3283
$perf.Range("B47").Value = "SCI"
$perf.Range("C47").Value = 24.1785
$perf.Range("B48").Value = "KICL"
$perf.Range("C48").Value = 24.1119
$perf.Range("B49").Value = "SKYGOLD"
$perf.Range("C49").Value = 23.9304
$perf.Range("B50").Value = "LORDSCHLO"
$perf.Range("C50").Value = 23.7541
$perf.Range("B51").Value = "AUBANK"
$perf.Range("C51").Value = 23.6403
$perf.Range("C52").Value = 23.2256
$perf.Range("B54").Value = "SURYODAY"
$perf.Range("C54").Value = 22.4892
$perf.Range("B55").Value = "INDIANB"
$perf.Range("C55").Value = 22.4463
$perf.Range("B56").Value = "GUJTHEM"
$perf.Range("C56").Value = 22.2937
$perf.Range("B58").Value = "ORBTEXP"
$perf.Range("C58").Value = 21.7352
$perf.Range("B59").Value = "TDPOWERSYS"
$perf.Range("C59").Value = 21.7288
$perf.Range("B60").Value = "CEATLTD"
$perf.Range("C60").Value = 20.1794
$perf.Range("B61").Value = "ATL"
$perf.Range("C61").Value = 20.1135
$perf.Range("B62").Value = "KAPSTON"
$perf.Range("C62").Value = 19.8943
$perf.Range("B63").Value = "USHAMART"
$perf.Range("C63").Value = 19.8759
$perf.Range("C64").Value = 19.6724
$perf.Range("B65").Value = "FEDERALBNK"
$perf.Range("C65").Value = 19.6311
$perf.Range("B66").Value = "GRMOVER"
$perf.Range("C66").Value = 19.3665
$perf.Range("B67").Value = "BANKINDIA"
$perf.Range("C67").Value = 19.2982
$perf.Range("C68").Value = 19.2187
$perf.Range("C70").Value = 19.0536
$perf.Range("C71").Value = 18.9965
$perf.Range("B73").Value = "MANAKCOAT"
$perf.Range("C73").Value = 18.8915
$perf.Range("C74").Value = 18.8394
$perf.Range("B75").Value = "THOMASCOTT"
$perf.Range("C75").Value = 18.6543
$perf.Range("B76").Value = "SHRIRAMFIN"
$perf.Range("C76").Value = 18.6082
